$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 259; existing rows 259:326 shift down to 260:327
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new record
$ws.Range("A259").Value = 3
$ws.Range("B259").Value = "Femacal de La Calera"
$ws.Range("C259").Value = "Coquimbo"
$ws.Range("D259").Value = 44722
$ws.Range("E259").Value = 5
$ws.Range("F259").Value = 100112039
$ws.Range("G259").Value = "Ciboulette"
$ws.Range("H259").Value = "Sin especificar"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 180
$ws.Range("K259").Value = 1500
$ws.Range("L259").Value = 1500
$ws.Range("M259").Value = 1500
$ws.Range("N259").Value = "$/docena de atados"
$ws.Range("O259").Value = "Provincia de Quillota"
$ws.Range("P259").Value = 500
$ws.Range("Q259").Value = 3
$ws.Range("R259").Value = "Hortaliza"
